$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SalesReturnHeader")
$ws1.Range("K2").Value = "Inclusive"
$ws1.Range("K2").Font.FontStyle = "Regular"
$ws1.Range("K2").Font.Name = "Calibri"
